$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9

# Row 12
$ws.Range("G12").Value = 2.65
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3.3
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 2.2
$ws.Range("W12").Value = 6.4
$ws.Range("X12").Value = 12
$ws.Range("Z12").Value = 32
$ws.Range("AA12").Value = 27
$ws.Range("AC12").Value = 5.8
$ws.Range("AD12").Value = 5.2
$ws.Range("AH12").Value = 7.3
$ws.Range("AI12").Value = 15
$ws.Range("AN12").Value = 4.35
$ws.Range("AO12").Value = 15
$ws.Range("AT12").Value = 2.18
$ws.Range("AU12").Value = 6.8
$ws.Range("AW12").Value = 4.8
$ws.Range("AY12").Value = 26
$ws.Range("AZ12").Value = 90
$ws.Range("BB12").Value = 350

# Row 14
$ws.Range("G14").Value = 1.75
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 4.75
$ws.Range("J14").Value = 2.4
$ws.Range("O14").Value = 1.36
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.15
$ws.Range("R14").Value = 1.67
$ws.Range("Y14").Value = 9
$ws.Range("AB14").Value = 34
$ws.Range("AH14").Value = 11
$ws.Range("AJ14").Value = 15
$ws.Range("AO14").Value = 9.5
$ws.Range("AQ14").Value = 34
$ws.Range("AX14").Value = 26
